$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quality_comparison")
$ws.Range("C2").Value = "approach"
Write-Host "Sheets:" $wb.Worksheets.Count
